# FB> editorial changes *
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the affiliation for Vanessa Didelez (row 13): was "University of Bremen"
$ws.Range("C13").Value = "Leibniz Institute for Prevention Research and Epidemiology - BIPS, Bremen"

# Correct the affiliation for Tracy Glass (row 23): was "University of Basel"
$ws.Range("C23").Value = "Swiss Tropical and Public Health Institute"

# Fix spelling of Fred Sorenson's surname (row 27): was "Sorensson"
$ws.Range("B27").Value = "Sorenson"

# Restore the active selection to B27 to match the saved view
$ws.Range("B27").Select()
